# Fruta / hortaliza, semanal
# A new week of "Brócoli" price data (Vega Monumental Concepción) is added.
# All existing weekly rows (168-177) shift down by one row (169-178), and
# the new top row (168) receives this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new last row (178) the same date-time number format used by the
# other rows in column D (style index reused automatically by the engine).
$ws.Range("D178").NumberFormat = $ws.Range("D177").NumberFormat

# Target state for rows 168-178, columns D (Fecha), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# O (Origen) and P (Precio $/Kg). All other columns (A,B,C,E,F,G,H,N,Q,R)
# are identical for every row in this block and remain unchanged.
$rows = @(
    @{ Row=168; D=44516; I="Primera"; J=1300; K=600; L=650; M=627; O="Región Metropolitana"; P=627 },
    @{ Row=169; D=44295; I="Primera"; J=1000; K=700; L=800; M=750; O="Región Metropolitana"; P=750 },
    @{ Row=170; D=44295; I="Segunda"; J=500;  K=600; L=600; M=600; O="Región Metropolitana"; P=600 },
    @{ Row=171; D=44509; I="Primera"; J=1000; K=700; L=800; M=750; O="Región Metropolitana"; P=750 },
    @{ Row=172; D=44509; I="Segunda"; J=500;  K=600; L=600; M=600; O="Región Metropolitana"; P=600 },
    @{ Row=173; D=44383; I="Primera"; J=1000; K=700; L=800; M=750; O="Región Metropolitana"; P=750 },
    @{ Row=174; D=44383; I="Segunda"; J=500;  K=600; L=600; M=600; O="Región Metropolitana"; P=600 },
    @{ Row=175; D=44273; I="Primera"; J=800;  K=800; L=900; M=850; O="Región del Maule";     P=850 },
    @{ Row=176; D=44273; I="Segunda"; J=400;  K=700; L=700; M=700; O="Región del Maule";     P=700 },
    @{ Row=177; D=44491; I="Primera"; J=1000; K=800; L=900; M=850; O="Región Metropolitana"; P=850 },
    @{ Row=178; D=44491; I="Segunda"; J=500;  K=700; L=700; M=700; O="Región Metropolitana"; P=700 }
)

foreach ($r in $rows) {
    $n = $r.Row

    # Columns that are constant across this whole block of rows; make sure
    # the brand-new row 178 also carries them (copy-safe no-op for the rest).
    $ws.Cells.Item($n, 1).Value2 = 11
    $ws.Cells.Item($n, 2).Value2 = "Vega Monumental Concepción"
    $ws.Cells.Item($n, 3).Value2 = "Bíobío"
    $ws.Cells.Item($n, 5).Value2 = 8
    $ws.Cells.Item($n, 6).Value2 = 100112023
    $ws.Cells.Item($n, 7).Value2 = "Brócoli"
    $ws.Cells.Item($n, 8).Value2 = "Sin especificar"
    $ws.Cells.Item($n, 14).Value2 = "`$/unidad"
    $ws.Cells.Item($n, 17).Value2 = 1
    $ws.Cells.Item($n, 18).Value2 = "Hortaliza"

    # Columns that actually change per the diff.
    $ws.Cells.Item($n, 4).Value2 = $r.D
    $ws.Cells.Item($n, 9).Value2 = $r.I
    $ws.Cells.Item($n, 10).Value2 = $r.J
    $ws.Cells.Item($n, 11).Value2 = $r.K
    $ws.Cells.Item($n, 12).Value2 = $r.L
    $ws.Cells.Item($n, 13).Value2 = $r.M
    $ws.Cells.Item($n, 15).Value2 = $r.O
    $ws.Cells.Item($n, 16).Value2 = $r.P
}
